# Auto-generated script applying value updates described by the diff
# (Ixion Profits workbook split into per-job sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 249.12
$ws.Range("J19").Value = 279.6154
$ws.Range("L19").Value = 279.6154
$ws.Range("N19").Value = -629.6154
$ws.Range("H62").Value = 1875.625
$ws.Range("I62").Value = 1972.1428
$ws.Range("J62").Value = 1200
$ws.Range("K62").Value = 1972.1428
$ws.Range("L62").Value = 1200
$ws.Range("M62").Value = -1348.1428
$ws.Range("N62").Value = -2448
$ws.Range("H65").Value = 1875.625
$ws.Range("I65").Value = 1972.1428
$ws.Range("J65").Value = 1200
$ws.Range("K65").Value = 9860.714
$ws.Range("L65").Value = 6000
$ws.Range("M65").Value = -6740.714
$ws.Range("N65").Value = -12240
$ws.Range("H111").Value = 101463.9
$ws.Range("I111").Value = 1071.6
$ws.Range("K111").Value = 3214.8
$ws.Range("M111").Value = -147.7999999999997
$ws.Range("H112").Value = 7813547
$ws.Range("I112").Value = 742.5
$ws.Range("J112").Value = 8334400.5
$ws.Range("K112").Value = 2227.5
$ws.Range("L112").Value = 25003201.5
$ws.Range("M112").Value = -1119.5
$ws.Range("N112").Value = -25005417.5
$ws.Range("H129").Value = 1012.8043
$ws.Range("J129").Value = 1182.3889
$ws.Range("L129").Value = 3547.1667
$ws.Range("N129").Value = -13547.1667
$ws.Range("H138").Value = 3141.875
$ws.Range("I138").Value = 1030.7646
$ws.Range("J138").Value = 5534.467
$ws.Range("K138").Value = 3092.2938
$ws.Range("L138").Value = 16603.401
$ws.Range("M138").Value = 2047.7062
$ws.Range("N138").Value = -26883.401

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2800.1667
$ws.Range("I2").Value = 3264
$ws.Range("K2").Value = 3264
$ws.Range("M2").Value = -3151
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H45").Value = 6011.5454
$ws.Range("I45").Value = 9978.637000000001
$ws.Range("K45").Value = 9978.637000000001
$ws.Range("M45").Value = -9601.637000000001
$ws.Range("H74").Value = 1695.4921
$ws.Range("I74").Value = 1559.037
$ws.Range("J74").Value = 2514.2222
$ws.Range("K74").Value = 1559.037
$ws.Range("L74").Value = 2514.2222
$ws.Range("M74").Value = -685.037
$ws.Range("N74").Value = -4262.2222
$ws.Range("H77").Value = 1695.4921
$ws.Range("I77").Value = 1559.037
$ws.Range("J77").Value = 2514.2222
$ws.Range("K77").Value = 7795.185
$ws.Range("L77").Value = 12571.111
$ws.Range("M77").Value = -3427.185
$ws.Range("N77").Value = -21307.111
$ws.Range("H102").Value = 12348145
$ws.Range("I102").Value = 12348145
$ws.Range("K102").Value = 12348145
$ws.Range("M102").Value = -12346523
$ws.Range("H116").Value = 2800.1667
$ws.Range("I116").Value = 3264
$ws.Range("K116").Value = 3264
$ws.Range("M116").Value = -970

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2800.1667
$ws.Range("I3").Value = 3264
$ws.Range("K3").Value = 3264
$ws.Range("M3").Value = -3150
$ws.Range("H20").Value = 1522.3125
$ws.Range("I20").Value = 1442
$ws.Range("J20").Value = 1656.1666
$ws.Range("K20").Value = 1442
$ws.Range("L20").Value = 1656.1666
$ws.Range("M20").Value = -1195
$ws.Range("N20").Value = -2150.1666
$ws.Range("H86").Value = 16668350
$ws.Range("I86").Value = 19609482
$ws.Range("J86").Value = 1933
$ws.Range("K86").Value = 19609482
$ws.Range("L86").Value = 1933
$ws.Range("M86").Value = -19608359
$ws.Range("N86").Value = -4179
$ws.Range("H89").Value = 16668350
$ws.Range("I89").Value = 19609482
$ws.Range("J89").Value = 1933
$ws.Range("K89").Value = 98047410
$ws.Range("L89").Value = 9665
$ws.Range("M89").Value = -98041794
$ws.Range("N89").Value = -20897
$ws.Range("H105").Value = 8686.8125
$ws.Range("I105").Value = 13085
$ws.Range("J105").Value = 3032
$ws.Range("K105").Value = 13085
$ws.Range("L105").Value = 3032
$ws.Range("M105").Value = -11338
$ws.Range("N105").Value = -6526

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6872.2607
$ws.Range("I31").Value = 819.0769
$ws.Range("J31").Value = 14741.4
$ws.Range("K31").Value = 819.0769
$ws.Range("L31").Value = 14741.4
$ws.Range("M31").Value = -524.0769
$ws.Range("N31").Value = -15331.4
$ws.Range("H34").Value = 6872.2607
$ws.Range("I34").Value = 819.0769
$ws.Range("J34").Value = 14741.4
$ws.Range("K34").Value = 819.0769
$ws.Range("L34").Value = 14741.4
$ws.Range("M34").Value = -617.0769
$ws.Range("N34").Value = -15145.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 310.45
$ws.Range("J107").Value = 408.41666
$ws.Range("L107").Value = 1225.24998
$ws.Range("N107").Value = -5065.249980000001
$ws.Range("H113").Value = 4286270.5
$ws.Range("I113").Value = 12500411
$ws.Range("J113").Value = 1000614.4
$ws.Range("K113").Value = 37501233
$ws.Range("L113").Value = 3001843.2
$ws.Range("M113").Value = -37499063
$ws.Range("N113").Value = -3006183.2
$ws.Range("H121").Value = 966.34784
$ws.Range("I121").Value = 400
$ws.Range("J121").Value = 1020.2857
$ws.Range("K121").Value = 1200
$ws.Range("L121").Value = 3060.8571
$ws.Range("M121").Value = 110
$ws.Range("N121").Value = -5680.8571
$ws.Range("H140").Value = 2263.9375
$ws.Range("I140").Value = 2263.9375
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 6791.8125
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -1611.8125
$ws.Range("N140").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1405.76
$ws.Range("I107").Value = 908.125
$ws.Range("J107").Value = 1639.9412
$ws.Range("K107").Value = 908.125
$ws.Range("L107").Value = 1639.9412
$ws.Range("M107").Value = 1011.875
$ws.Range("N107").Value = -5479.9412
$ws.Range("H126").Value = 10292.75
$ws.Range("I126").Value = 12879.111
$ws.Range("J126").Value = 2533.6667
$ws.Range("K126").Value = 38637.333
$ws.Range("L126").Value = 7601.000100000001
$ws.Range("M126").Value = -36167.333
$ws.Range("N126").Value = -12541.0001
$ws.Range("H132").Value = 2782.697
$ws.Range("I132").Value = 2495.054
$ws.Range("J132").Value = 3149.6897
$ws.Range("K132").Value = 7485.162
$ws.Range("L132").Value = 9449.069100000001
$ws.Range("M132").Value = -4955.162
$ws.Range("N132").Value = -14509.0691
$ws.Range("H141").Value = 111000
$ws.Range("J141").Value = 111000
$ws.Range("L141").Value = 111000
$ws.Range("N141").Value = -121360

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3019.8
$ws.Range("I7").Value = 2024.75
$ws.Range("J7").Value = 7000
$ws.Range("K7").Value = 2024.75
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = -1912.75
$ws.Range("N7").Value = -7224
$ws.Range("H40").Value = 111113690
$ws.Range("I40").Value = 125001780
$ws.Range("K40").Value = 125001780
$ws.Range("M40").Value = -125001644
$ws.Range("H55").Value = 107143304
$ws.Range("I55").Value = 200000240
$ws.Range("J55").Value = 55556110
$ws.Range("K55").Value = 200000240
$ws.Range("L55").Value = 55556110
$ws.Range("M55").Value = -200000067
$ws.Range("N55").Value = -55556456
$ws.Range("H122").Value = 5497183
$ws.Range("I122").Value = 5955031.5
$ws.Range("K122").Value = 17865094.5
$ws.Range("M122").Value = -17862644.5
$ws.Range("H126").Value = 3019.8
$ws.Range("I126").Value = 2024.75
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 6074.25
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -3604.25
$ws.Range("N126").Value = -25940
$ws.Range("H132").Value = 16673529
$ws.Range("I132").Value = 21674388
$ws.Range("J132").Value = 4000.6667
$ws.Range("K132").Value = 65023164
$ws.Range("L132").Value = 12002.0001
$ws.Range("M132").Value = -65020634
$ws.Range("N132").Value = -17062.0001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 47619430
$ws.Range("I107").Value = 71428950
$ws.Range("J107").Value = 388.42856
$ws.Range("K107").Value = 214286850
$ws.Range("L107").Value = 1165.28568
$ws.Range("M107").Value = -214284930
$ws.Range("N107").Value = -5005.28568
$ws.Range("H126").Value = 2026.091
$ws.Range("I126").Value = 1335.75
$ws.Range("J126").Value = 2420.5715
$ws.Range("K126").Value = 4007.25
$ws.Range("L126").Value = 7261.7145
$ws.Range("M126").Value = -1537.25
$ws.Range("N126").Value = -12201.7145
$ws.Range("H135").Value = 45083
$ws.Range("J135").Value = 45083
$ws.Range("L135").Value = 45083
$ws.Range("N135").Value = -55223
$ws.Range("H136").Value = 1445.9524
$ws.Range("I136").Value = 840.9091
$ws.Range("J136").Value = 2111.5
$ws.Range("K136").Value = 2522.7273
$ws.Range("L136").Value = 6334.5
$ws.Range("M136").Value = 27.27269999999999
$ws.Range("N136").Value = -11434.5
